$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny floating-point precision refinements to existing "Adj Close" (F) values ---
# --- plus corrected High (C214) and Volume (G214) for the last existing row      ---
$ws.Range("F2").Value = 57.429565
$ws.Range("F4").Value = 57.595856
$ws.Range("F5").Value = 58.153419
$ws.Range("F6").Value = 58.261024
$ws.Range("F10").Value = 58.46645
$ws.Range("F11").Value = 58.427311
$ws.Range("F13").Value = 58.838158
$ws.Range("F15").Value = 58.446884
$ws.Range("F18").Value = 58.554482
$ws.Range("F19").Value = 58.46645
$ws.Range("F20").Value = 58.339287
$ws.Range("F21").Value = 58.564266
$ws.Range("F23").Value = 58.88969
$ws.Range("F24").Value = 58.399597
$ws.Range("F26").Value = 57.919296
$ws.Range("F27").Value = 58.056526
$ws.Range("F28").Value = 57.870293
$ws.Range("F32").Value = 57.429214
$ws.Range("F33").Value = 57.340988
$ws.Range("F35").Value = 57.04694
$ws.Range("F37").Value = 57.24297
$ws.Range("F38").Value = 56.968517
$ws.Range("F40").Value = 57.135155
$ws.Range("F42").Value = 56.68829
$ws.Range("F43").Value = 56.943638
$ws.Range("F45").Value = 56.796322
$ws.Range("F46").Value = 56.747211
$ws.Range("F47").Value = 57.100784
$ws.Range("F48").Value = 57.837368
$ws.Range("F49").Value = 58.534683
$ws.Range("F50").Value = 58.132011
$ws.Range("F52").Value = 58.397182
$ws.Range("F53").Value = 58.966816
$ws.Range("F54").Value = 58.760574
$ws.Range("F55").Value = 58.348083
$ws.Range("F56").Value = 58.956993
$ws.Range("F57").Value = 59.310562
$ws.Range("F58").Value = 59.29092
$ws.Range("F60").Value = 58.652534
$ws.Range("F61").Value = 58.573971
$ws.Range("F62").Value = 58.642712
$ws.Range("F63").Value = 58.888248
$ws.Range("F65").Value = 59.414783
$ws.Range("F66").Value = 59.562408
$ws.Range("F67").Value = 59.572254
$ws.Range("F70").Value = 59.277
$ws.Range("F73").Value = 58.656967
$ws.Range("F74").Value = 58.716019
$ws.Range("F77").Value = 58.784904
$ws.Range("F80").Value = 59.326206
$ws.Range("F85").Value = 59.752251
$ws.Range("F88").Value = 59.278893
$ws.Range("F89").Value = 59.229584
$ws.Range("F91").Value = 59.70295
$ws.Range("F93").Value = 59.387371
$ws.Range("F95").Value = 59.061932
$ws.Range("F96").Value = 58.785809
$ws.Range("F97").Value = 58.598431
$ws.Range("F98").Value = 58.558987
$ws.Range("F101").Value = 58.174374
$ws.Range("F103").Value = 58.529404
$ws.Range("F104").Value = 58.687183
$ws.Range("F107").Value = 58.458889
$ws.Range("F109").Value = 58.192039
$ws.Range("F110").Value = 58.458889
$ws.Range("F111").Value = 58.271103
$ws.Range("F113").Value = 58.063557
$ws.Range("F114").Value = 58.05368
$ws.Range("F115").Value = 58.409466
$ws.Range("F116").Value = 58.192039
$ws.Range("F117").Value = 58.290867
$ws.Range("F118").Value = 58.320518
$ws.Range("F120").Value = 58.231575
$ws.Range("F124").Value = 57.925198
$ws.Range("F128").Value = 57.318913
$ws.Range("F129").Value = 57.318913
$ws.Range("F130").Value = 57.586349
$ws.Range("F131").Value = 57.625961
$ws.Range("F132").Value = 58.071678
$ws.Range("F133").Value = 58.487679
$ws.Range("F134").Value = 58.220245
$ws.Range("F138").Value = 58.1311
$ws.Range("F144").Value = 57.873581
$ws.Range("F147").Value = 57.70274
$ws.Range("F149").Value = 57.911194
$ws.Range("F150").Value = 57.841713
$ws.Range("F154").Value = 57.504208
$ws.Range("F155").Value = 57.39502
$ws.Range("F161").Value = 57.067444
$ws.Range("F163").Value = 57.345383
$ws.Range("F165").Value = 57.385094
$ws.Range("F167").Value = 57.70274
$ws.Range("F168").Value = 57.831783
$ws.Range("C214").Value = 57.18
$ws.Range("G214").Value = 2679700

# --- Append two new trading days (rows 215-216) ---
# Column A holds dates as plain text (matches the rest of the sheet), so force
# text formatting before assignment, then clear the format stamp it leaves behind.
$ws.Range("A215").NumberFormat = "@"
$ws.Range("A215").Value = "2023-11-07"
$ws.Range("A215").ClearFormats()
$ws.Range("B215").Value = 57.18
$ws.Range("C215").Value = 57.360001
$ws.Range("D215").Value = 57.169998
$ws.Range("E215").Value = 57.27
$ws.Range("F215").Value = 57.27
$ws.Range("G215").Value = 2038000
$ws.Range("A216").NumberFormat = "@"
$ws.Range("A216").Value = "2023-11-08"
$ws.Range("A216").ClearFormats()
$ws.Range("B216").Value = 57.259998
$ws.Range("C216").Value = 57.41
$ws.Range("D216").Value = 57.25
$ws.Range("E216").Value = 57.369999
$ws.Range("F216").Value = 57.369999
$ws.Range("G216").Value = 2274015
